$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ticker table for columns B:F, rows 2-50 (49 data rows).
# Row index 0 of this array corresponds to worksheet row 2.
$data = @(
    @('NSE:AARTIDRUGS','NSE:ESSENTIA','NSE:COROMANDEL','NSE:GAIL','NSE:ASIANPAINT'),
    @('NSE:AGRITECH','NSE:LINDEINDIA','','NSE:POWERGRID','NSE:BAJAJ-AUTO'),
    @('NSE:ALKYLAMINE','NSE:ORCHPHARMA','','','NSE:COROMANDEL'),
    @('NSE:ALPHAETF','NSE:ORISSAMINE','','','NSE:GRASIM'),
    @('NSE:APTUS','','','','NSE:HINDUNILVR'),
    @('NSE:ARVINDFASN','','','','NSE:M&MFIN'),
    @('NSE:ASIANPAINT','','','',''),
    @('NSE:AVROIND','','','',''),
    @('NSE:BAJAJ-AUTO','','','',''),
    @('NSE:BALAXI','','','',''),
    @('NSE:BALPHARMA','','','',''),
    @('NSE:BANSWRAS','','','',''),
    @('NSE:COASTCORP','','','',''),
    @('NSE:DABUR','','','',''),
    @('NSE:DEEPENR','','','',''),
    @('NSE:DHANBANK','','','',''),
    @('NSE:ELECTCAST','','','',''),
    @('NSE:ERIS','','','',''),
    @('NSE:ESTER','','','',''),
    @('NSE:GKWLIMITED','','','',''),
    @('NSE:GRASIM','','','',''),
    @('NSE:GSFC','','','',''),
    @('NSE:GSPL','','','',''),
    @('NSE:GUFICBIO','','','',''),
    @('NSE:GULFOILLUB','','','',''),
    @('NSE:GULPOLY','','','',''),
    @('NSE:HDFCMID150','','','',''),
    @('NSE:HINDUNILVR','','','',''),
    @('NSE:HIRECT','','','',''),
    @('NSE:INDIGOPNTS','','','',''),
    @('NSE:INOXGREEN','','','',''),
    @('NSE:IPL','','','',''),
    @('NSE:ITI','','','',''),
    @('NSE:JAYSREETEA','','','',''),
    @('NSE:KANPRPLA','','','',''),
    @('NSE:LUMAXIND','','','',''),
    @('NSE:LYPSAGEMS','','','',''),
    @('NSE:MANAPPURAM','','','',''),
    @('NSE:MIRZAINT','','','',''),
    @('NSE:MUTHOOTCAP','','','',''),
    @('NSE:NFL','','','',''),
    @('NSE:ORIENTALTL','','','',''),
    @('NSE:PATINTLOG','','','',''),
    @('NSE:PRECAM','','','',''),
    @('NSE:PRITI','','','',''),
    @('NSE:RCF','','','',''),
    @('NSE:RELIGARE','','','',''),
    @('NSE:RUBYMILLS','','','',''),
    @('NSE:RUCHINFRA','','','','')
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt 5; $j++) {
        $col = $j + 2
        $ws.Cells.Item($r, $col).Value = $rowVals[$j]
    }
}

# The sheet previously extended to row 62; the refreshed list only needs
# rows 2-50, so remove the now-unused trailing rows entirely.
$ws.Range("A51:F62").EntireRow.Delete()
